$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.120.43'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.792.30'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.63'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5230'
$ws.Range('E7').Value = '  +3.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3796'
$ws.Range('E8').Value = '  -3.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07954'
$ws.Range('E9').Value = '  -3.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.43'
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.092'
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.246'
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.45'
$ws.Range('E14').Value = '  -2.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.797.44'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.281'
$ws.Range('E16').Value = '  -2.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.53'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001086'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06575'
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.28'
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.956'
$ws.Range('E22').Value = '  -2.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.166.93'
$ws.Range('E23').Value = '  -0.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.10'
$ws.Range('E24').Value = '  -1.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.259'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.28'
$ws.Range('E26').Value = '  +3.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.38'
$ws.Range('E27').Value = '  -3.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.998.32'
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.332'
$ws.Range('E29').Value = '  -2.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.42'
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1084'
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.052'
$ws.Range('E32').Value = '  -4.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.691'
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.519'
$ws.Range('E34').Value = '  -4.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07196'
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.19'
$ws.Range('E36').Value = '  +8.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02304'
$ws.Range('E37').Value = '  -1.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2143'
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('E39').Value = '  -3.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.584'
$ws.Range('E40').Value = '  -2.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6148'
$ws.Range('E41').Value = '  -1.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.161'
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.374'
$ws.Range('E43').Value = '  -2.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.15'
$ws.Range('E44').Value = '  -2.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.770'
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5920'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '126.49'
$ws.Range('E47').Value = '  +1.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.214'
$ws.Range('E48').Value = '  +2.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.916'
$ws.Range('E49').Value = '  -2.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06780'
$ws.Range('E50').Value = '  -1.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.54'
$ws.Range('E51').Value = '  -1.93%  '
